$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the shared formula in column B (rows 1 to 77): change ".dump" -> ".s"
for ($r = 1; $r -le 77; $r++) {
    $cellA = "A" + $r
    $cellB = "B" + $r
    $ws.Range($cellB).Formula = '=CONCATENATE("riscv32-unknown-elf-objdump -d test/", ' + $cellA + ', ".riscv"," >  test/", ' + $cellA + ', ".s")'
}

# Set the selection on the sheet view: active cell B1, selected range B1:B77
$ws.Range("B1:B77").Select()
